$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns: F (president) and H (senate_seats).
# Old F (senate_party_in_power) shifts to G; old G (house_party_in_power) shifts to I.
$ws.Columns("F:F").Insert()
$ws.Columns("H:H").Insert()

# Header row
$ws.Range("F1").Value = "president"
$ws.Range("H1").Value = "senate_seats"
$ws.Range("J1").Value = "house_seats"

# Data rows
$ws.Range("F2").Value = "Ford"
$ws.Range("H2").Value = "62 - 38"
$ws.Range("J2").Value = "291 - 144"
$ws.Range("F3").Value = "Carter"
$ws.Range("H3").Value = "59 - 41"
$ws.Range("J3").Value = "292 - 143"
$ws.Range("F4").Value = "Carter"
$ws.Range("H4").Value = "59 - 41"
$ws.Range("J4").Value = "292 - 143"
$ws.Range("F5").Value = "Carter"
$ws.Range("H5").Value = "59 - 41"
$ws.Range("J5").Value = "292 - 143"
$ws.Range("F6").Value = "Carter"
$ws.Range("H6").Value = "59 - 41"
$ws.Range("J6").Value = "292 - 143"
$ws.Range("F7").Value = "Carter"
$ws.Range("H7").Value = "58 - 42"
$ws.Range("J7").Value = "277 - 158"
$ws.Range("F8").Value = "Reagan"
$ws.Range("H8").Value = "53 - 47"
$ws.Range("J8").Value = "244 - 191"
$ws.Range("F9").Value = "Reagan"
$ws.Range("H9").Value = "53 - 47"
$ws.Range("J9").Value = "244 - 191"
$ws.Range("F10").Value = "Reagan"
$ws.Range("H10").Value = "53 - 47"
$ws.Range("J10").Value = "244 - 191"
$ws.Range("F11").Value = "Reagan"
$ws.Range("H11").Value = "55 - 45"
$ws.Range("J11").Value = "271 - 164"
$ws.Range("F12").Value = "Reagan"
$ws.Range("H12").Value = "55 - 45"
$ws.Range("J12").Value = "270 - 165"
$ws.Range("F13").Value = "Reagan"
$ws.Range("H13").Value = "55 - 45"
$ws.Range("J13").Value = "270 - 165"
$ws.Range("F14").Value = "Reagan"
$ws.Range("H14").Value = "53 - 47"
$ws.Range("J14").Value = "253 - 182"
$ws.Range("F15").Value = "Reagan"
$ws.Range("H15").Value = "54 - 46"
$ws.Range("J15").Value = "258 - 177"
$ws.Range("F16").Value = "G. H. W. Bush"
$ws.Range("H16").Value = "55 - 45"
$ws.Range("J16").Value = "258 - 176"
$ws.Range("F17").Value = "Clinton"
$ws.Range("H17").Value = "53 - 47"
$ws.Range("J17").Value = "233 - 199"
$ws.Range("F18").Value = "Clinton"
$ws.Range("H18").Value = "53 - 47"
$ws.Range("J18").Value = "235 - 198"
$ws.Range("F19").Value = "Obama"
$ws.Range("H19").Value = "54 - 46"
$ws.Range("J19").Value = "232 - 200"
$ws.Range("F20").Value = "Trump"
$ws.Range("H20").Value = "51 - 49"
$ws.Range("J20").Value = "238 - 193"
$ws.Range("F21").Value = "Trump"
$ws.Range("H21").Value = "51 - 49"
$ws.Range("J21").Value = "238 - 193"

# Column J is brand new (no neighbor to inherit formatting from via the
# column-insert shift), so match the centered style used by the other
# data columns (E:I) explicitly.
$ws.Range("J2:J21").HorizontalAlignment = -4108

# Autofit the new columns to size them like the existing data columns
$ws.Columns("F:F").AutoFit()
$ws.Columns("H:H").AutoFit()
$ws.Columns("J:J").AutoFit()

# Match the author's final selection
$ws.Range("A2").Select()
